$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 0.6226591760299626
$summary.Range("C2").Value = 0.5771495877502945
$summary.Range("D2").Value = 0.9176029962546817
$summary.Range("E2").Value = 0.7086044830079538
$summary.Range("F2").Value = 0.8207705192629816
$summary.Range("G2").Value = 0.8972462849496443
$summary.Range("H2").Value = 0.7754176661195977
$summary.Range("I2").Value = 490
$summary.Range("J2").Value = 359
$summary.Range("K2").Value = 175
$summary.Range("L2").Value = 44

# --- Sheet 2: Classification Report ---
$classRep = $wb.Worksheets.Item("Classification Report")

# Row 2 ("0")
$classRep.Range("B2").Value = 0.7990867579908676
$classRep.Range("C2").Value = 0.3277153558052435
$classRep.Range("D2").Value = 0.4648074369189907

# Row 3 ("1")
$classRep.Range("B3").Value = 0.5771495877502945
$classRep.Range("C3").Value = 0.9176029962546817
$classRep.Range("D3").Value = 0.7086044830079538

# Row 4 ("accuracy")
$classRep.Range("B4").Value = 0.6226591760299626
$classRep.Range("C4").Value = 0.6226591760299626
$classRep.Range("D4").Value = 0.6226591760299626
$classRep.Range("E4").Value = 0.6226591760299626

# Row 5 ("macro avg")
$classRep.Range("B5").Value = 0.688118172870581
$classRep.Range("C5").Value = 0.6226591760299626
$classRep.Range("D5").Value = 0.5867059599634722

# Row 6 ("weighted avg")
$classRep.Range("B6").Value = 0.6881181728705811
$classRep.Range("C6").Value = 0.6226591760299626
$classRep.Range("D6").Value = 0.5867059599634722

# --- Sheet 3: Confusion Matrix ---
$confMat = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 ("Actual 0")
$confMat.Range("B2").Value = 175
$confMat.Range("C2").Value = 359

# Row 3 ("Actual 1")
$confMat.Range("B3").Value = 44
$confMat.Range("C3").Value = 490
